# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect refreshed scrape data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - first data table
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 75
$wsExpo.Range("F7").Value  = 592
$wsExpo.Range("F8").Value  = 115
$wsExpo.Range("F9").Value  = 8765
$wsExpo.Range("F13").Value = 990
$wsExpo.Range("F14").Value = 113
$wsExpo.Range("F18").Value = 265
$wsExpo.Range("F20").Value = 232
$wsExpo.Range("F21").Value = 1043

# Sheet "全部类型" (all types) - combined table with the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 75
$wsAll.Range("F9").Value  = 592
$wsAll.Range("F10").Value = 115
$wsAll.Range("F11").Value = 8765
$wsAll.Range("F15").Value = 990
$wsAll.Range("F16").Value = 113
$wsAll.Range("F20").Value = 265
$wsAll.Range("F22").Value = 232
$wsAll.Range("F23").Value = 1043
